$d = $word.ActiveDocument

$newText = "Сазвежђе Perseus: 16. до 25. јануара, 7. и 16. новембра, 6. до 15. децембра"

# The first occurrence of this paragraph is wrapped in a leftover Word
# "_Hlk..." auto-bookmark (from an old hyperlink) that must be removed
# entirely along with the text runs it spans.
try {
    $bm = $d.Bookmarks.Item("_Hlk514861060")
    $bm.Delete()
} catch {
    # bookmark already absent - nothing to do
}

# Every paragraph whose full text is the old (multi-run) Serbian sentence
# about the Perseus constellation observation window gets collapsed down
# to a single, unformatted run containing the new translated sentence.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*2018*8. децембра*") {
        # Drop every existing run (and any run-level formatting) in the
        # paragraph, but keep the paragraph mark itself (and thus the
        # paragraph's own pPr / sectPr) intact.
        $r = $p.Range
        [void]$r.MoveEnd(1, -1)
        $r.Delete()

        # Re-fetch the now-empty paragraph and insert the single plain
        # (unformatted) replacement run before the paragraph mark.
        $p2 = $d.Paragraphs.Item($i)
        $r2 = $p2.Range
        [void]$r2.MoveEnd(1, -1)
        $r2.InsertAfter($newText)
    }
}
